$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row above row 29, shifting existing rows 29-53 down to 30-54
$ws.Rows.Item(29).Insert()

# Populate the new row 29 with the September details/date values
$ws.Cells.Item(29, 18).Value = "balance your axis"
$ws.Cells.Item(29, 19).Value = "2024-09-04 08:14:16"
